$p = $ppt.ActivePresentation

# 1) Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") has a 2-column table
#    (the graphicFrame shape on that slide) whose table style is switched
#    to a different built-in style id.
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{49C19E07-75CC-4B3E-8802-5C050A47F215}")
    }
}

# 2) The deck's theme colour scheme is swapped from the "Integral" palette
#    to the default "Office" palette (dk1/lt1 - black/white - stay the
#    same; the other ten slots change).
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
